$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phase 1: product codes (A) and names (B) for all new rows, in row order
$ws.Range("A2892").Value = "TM.13156"
$ws.Range("B2892").Value = "ANDELİ STABLİZATOR MOTORLU 500W AD41 Q/4 P/1 X"
$ws.Range("A2893").Value = "TM.13157"
$ws.Range("B2893").Value = "STABLİZATOR ANDELİ MOTORLU 1000W 150-250V AD42 Q/4 P/1"
$ws.Range("A2894").Value = "TM.14648"
$ws.Range("B2894").Value = "STABLİZATOR ANDELİ MOTORLU 7500W AD49 VERTIKAL X"
$ws.Range("A2895").Value = "TM.14649"
$ws.Range("B2895").Value = "STABLİZATOR ANDELİ MOTORLU 10000W AD50VERTIKALQ/1X"
$ws.Range("A2896").Value = "TM.14650"
$ws.Range("B2896").Value = "ANDELİ QAYNAQ APARATI ARC 160G AD51 Q/1 X"
$ws.Range("A2897").Value = "TM.14651"
$ws.Range("B2897").Value = "ANDELİ QAYNAQ APARATI ARC 180G AD52 Q/1 X"
$ws.Range("A2898").Value = "TM.14653"
$ws.Range("B2898").Value = "QAYNAQ APARATI ANDELİ TIG 180G AD54 Q/1 (M) X"
$ws.Range("A2899").Value = "TM.14654"
$ws.Range("B2899").Value = "ANDELİ QAYNAQ APARATI TIG 200G AD55 Q/1 X"
$ws.Range("A2900").Value = "TM.14655"
$ws.Range("B2900").Value = "QAYNAQ APARATI ANDELİ  TIG 160 AD56 Q/1 X"
$ws.Range("A2901").Value = "TM.14656"
$ws.Range("B2901").Value = "QAYNAQ APARATI ANDELİ TIG 180 AD57 Q/1 X"
$ws.Range("A2902").Value = "TM.14664"
$ws.Range("B2902").Value = "ANDELİ ELEKT.AVTOMATI MCCB,AM1-400L/3P 250A AD66 X"
$ws.Range("A2903").Value = "TM.14666"
$ws.Range("B2903").Value = "ANDELİ ELEKT.AVTOMATI MCCB,AM3T-160/3P 160A AD68 X"
$ws.Range("A2904").Value = "TM.14667"
$ws.Range("B2904").Value = "ANDELİ ELEKT.AVTOMATI MCCB,AM3T-250/3P 250A AD69 X"
$ws.Range("A2905").Value = "TM.14668"
$ws.Range("B2905").Value = "ANDELİ ELEKT.AVTOMAT MCCB,AM8-160/3P 160A AD70Q/1X"
$ws.Range("A2906").Value = "TM.14669"
$ws.Range("B2906").Value = "ANDELİ ELEKT.AVTOMAT MCCB,AM8-250/3P 250A AD71Q/1X"
$ws.Range("A2907").Value = "TM.15800"
$ws.Range("B2907").Value = "ANDELİ QAYNAQ APARATI MIG200Y AD73 Q/1 X"
$ws.Range("A2908").Value = "TM.15801"
$ws.Range("B2908").Value = "ANDELİ QAYNAQ APARATI MIG250Y AD74 Q/1 X"
$ws.Range("A2909").Value = "TM.15802"
$ws.Range("B2909").Value = "ANDELİ QAYNAQ APARATI MIG350 AD75 Q/1 X"
$ws.Range("A2910").Value = "TM.15803"
$ws.Range("B2910").Value = "ANDELİ QAYNAQ APARATI ARC160 AD76 Q/1"
$ws.Range("A2911").Value = "TM.15804"
$ws.Range("B2911").Value = "ANDELİ QAYNAQ APARATI ARC200-S180A AD78 Q/1"
$ws.Range("A2912").Value = "TM.15805"
$ws.Range("B2912").Value = "ANDELİ QAYNAQ APARATI ARC200-S200A AD79 Q/1"
$ws.Range("A2913").Value = "TM.18150"
$ws.Range("B2913").Value = "STABLİZATOR ANDELİ ABVR-500VA (KOMBİ) AD90 Q/4 X"
$ws.Range("A2914").Value = "TM.18151"
$ws.Range("B2914").Value = "STABLİZATOR ANDELİ SXEMA KOMBİ Q/100 X"
$ws.Range("A2915").Value = "TM.220508010556"
$ws.Range("B2915").Value = "VOLTAJ TƏNZİMLƏYİCİ ANDELİ 63A/2P AD100 Q/1"
$ws.Range("A2916").Value = "TM.230508010797"
$ws.Range("B2916").Value = "STABLİZATOR ANDELİ SVC-40KVA-3 AD340 Q/1"
$ws.Range("A2917").Value = "TM.230508010798"
$ws.Range("B2917").Value = "STABLİZATOR ANDELİ SVC-50KVA-3 AD350 Q/1"
$ws.Range("A2918").Value = "TM.15790"
$ws.Range("B2918").Value = "ANDELİ STABLİZATOR DİGİTAL 500W AD80 Q/2 X"
$ws.Range("A2919").Value = "TM.15798"
$ws.Range("B2919").Value = "QAYNAQ APARATI ANDELİ ARC 315G AD80 Q/1 X"
$ws.Range("A2920").Value = "TM.15799"
$ws.Range("B2920").Value = "QAYNAQ APARATI ANDELİ TIG 250G AD72 Q/1 X"
$ws.Range("A2921").Value = "TM.200508010775"
$ws.Range("B2921").Value = "STABLİZATOR ANDELİ MOTORLU 5000W 150-250V AD46 Q/1 P/1"
$ws.Range("A2922").Value = "TM.200508010776"
$ws.Range("B2922").Value = "STABLİZATOR ANDELİ MOTORLU 7500W 150-250V AD47 Q/1 P/1"
$ws.Range("A2923").Value = "TM.200508010777"
$ws.Range("B2923").Value = "STABLİZATOR ANDELİ MOTORLU 10000W 150-250V AD48 Q/1 P/1"
$ws.Range("A2924").Value = "TM.200508010930"
$ws.Range("B2924").Value = "STABLİZATOR ANDELİ MOTORLU SVC-15000W 150-250V AD88 Q/1"
$ws.Range("A2925").Value = "TM.200508010931"
$ws.Range("B2925").Value = "STABLİZATOR ANDELİ MOTORLU SVC-20000W 150-250V AD89 Q/1"
$ws.Range("A2926").Value = "TM.200508011789"
$ws.Range("B2926").Value = "STABLİZATOR ANDELİ MOTORLU SVC-30000W 150-250V AD91 Q/1"
$ws.Range("A2927").Value = "TM.200508012078"
$ws.Range("B2927").Value = "STABLİZATOR ANDELİ SVC-3KVA 110-250V AD95 Q/1"
$ws.Range("A2928").Value = "TM.200508010772"
$ws.Range("B2928").Value = "STABLİZATOR ANDELİ MOTORLU 1500W 150-250V AD43 Q/4 P/1"
$ws.Range("A2929").Value = "TM.200508010773"
$ws.Range("B2929").Value = "STABLİZATOR ANDELİ MOTORLU 2000W 150-250V AD44 Q/1 P/1"
$ws.Range("A2930").Value = "TM.200508010774"
$ws.Range("B2930").Value = "STABLİZATOR ANDELİ MOTORLU 3000W 150-250V AD45 Q/1 P/1"
$ws.Range("A2931").Value = "TM.200508012079"
$ws.Range("B2931").Value = "STABLİZATOR ANDELİ SVC-5KVA 110-250V AD96 Q/1"
$ws.Range("A2932").Value = "TM.200508013197"
$ws.Range("B2932").Value = "STABLİZATOR ANDELİ SVC-7.5KVA 110-250V AD97 Q/1"
$ws.Range("A2933").Value = "TM.200508013198"
$ws.Range("B2933").Value = "STABLİZATOR ANDELİ SVC-10KVA 110-250V AD98 Q/1"
$ws.Range("A2934").Value = "TM.200508013473"
$ws.Range("B2934").Value = "STABLİZATOR ANDELİ SVC-15KVA 110-250V AD99 Q/1"
$ws.Range("A2935").Value = "TM.240508011773"
$ws.Range("B2935").Value = "STABLİZATOR ANDELİ SVC-20KVA 110-250V AD81 Q/1"
$ws.Range("A2936").Value = "TM.240508011774"
$ws.Range("B2936").Value = "STABLİZATOR ANDELİ SVC-30KVA 110-250V AD82 Q/1"
$ws.Range("A2937").Value = "TM.240805012256"
$ws.Range("B2937").Value = "STABLİZATOR ANDELİ AD110-40 SVC-40KVA 110-250V Q/1"
$ws.Range("A2938").Value = "TM.240805012257"
$ws.Range("B2938").Value = "STABLİZATOR ANDELİ AD110-50 SVC-50KVA 110-250V Q/1"
$ws.Range("A2939").Value = "TM.14657"
$ws.Range("B2939").Value = "ANDELİ ELEKT.AVTOMATI MCCB,AM1-63L/3P 50A AD58 X"
$ws.Range("A2940").Value = "TM.14658"
$ws.Range("B2940").Value = "ANDELİ ELEKT.AVTOMATI MCCB,AM1-63L/3P 63A AD59 X"
$ws.Range("A2941").Value = "TM.14659"
$ws.Range("B2941").Value = "ANDELİ ELEKT.AVTOMATI MCCB,AM-100L/3P80A AD61 Q/1X"
$ws.Range("A2942").Value = "TM.14661"
$ws.Range("B2942").Value = "ANDELİ ELEKT.AVTOMATI MCCB,AM-225L/3P125A AD63Q/1X"
$ws.Range("A2943").Value = "TM.14662"
$ws.Range("B2943").Value = "ANDELİ ELEKT.AVTOMATI MCCB,AM-225L/3P 160A AD64 X"

# Phase 2: price (C) and classification columns C-G for all new rows, in row order
$ws.Range("C2892").Value = 38.4
$ws.Range("D2892").Value = "Elektrik"
$ws.Range("E2892").Value = "Elektrik"
$ws.Range("F2892").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2892").Value = "Stablizatorlar"
$ws.Range("C2893").Value = 44.2
$ws.Range("D2893").Value = "Elektrik"
$ws.Range("E2893").Value = "Elektrik"
$ws.Range("F2893").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2893").Value = "Stablizatorlar"
$ws.Range("C2894").Value = 295
$ws.Range("D2894").Value = "Elektrik"
$ws.Range("E2894").Value = "Elektrik"
$ws.Range("F2894").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2894").Value = "Stablizatorlar"
$ws.Range("C2895").Value = 372
$ws.Range("D2895").Value = "Elektrik"
$ws.Range("E2895").Value = "Elektrik"
$ws.Range("F2895").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2895").Value = "Stablizatorlar"
$ws.Range("C2896").Value = 149
$ws.Range("D2896").Value = "Elektrik"
$ws.Range("E2896").Value = "Elektrikli Əl Alətləri"
$ws.Range("F2896").Value = "Qaynaq Aparatları və Aksesuarları"
$ws.Range("G2896").Value = "Metal Qaynaq Aparatları və Aksesuarları"
$ws.Range("C2897").Value = 278
$ws.Range("D2897").Value = "Elektrik"
$ws.Range("E2897").Value = "Elektrikli Əl Alətləri"
$ws.Range("F2897").Value = "Qaynaq Aparatları və Aksesuarları"
$ws.Range("G2897").Value = "Metal Qaynaq Aparatları və Aksesuarları"
$ws.Range("C2898").Value = 195
$ws.Range("D2898").Value = "Elektrik"
$ws.Range("E2898").Value = "Elektrikli Əl Alətləri"
$ws.Range("F2898").Value = "Qaynaq Aparatları və Aksesuarları"
$ws.Range("G2898").Value = "Metal Qaynaq Aparatları və Aksesuarları"
$ws.Range("C2899").Value = 231
$ws.Range("D2899").Value = "Elektrik"
$ws.Range("E2899").Value = "Elektrikli Əl Alətləri"
$ws.Range("F2899").Value = "Qaynaq Aparatları və Aksesuarları"
$ws.Range("G2899").Value = "Metal Qaynaq Aparatları və Aksesuarları"
$ws.Range("C2900").Value = 330
$ws.Range("D2900").Value = "Elektrik"
$ws.Range("E2900").Value = "Elektrikli Əl Alətləri"
$ws.Range("F2900").Value = "Qaynaq Aparatları və Aksesuarları"
$ws.Range("G2900").Value = "Metal Qaynaq Aparatları və Aksesuarları"
$ws.Range("C2901").Value = 387
$ws.Range("D2901").Value = "Elektrik"
$ws.Range("E2901").Value = "Elektrikli Əl Alətləri"
$ws.Range("F2901").Value = "Qaynaq Aparatları və Aksesuarları"
$ws.Range("G2901").Value = "Metal Qaynaq Aparatları və Aksesuarları"
$ws.Range("C2902").Value = 110.97
$ws.Range("D2902").Value = "Elektrik"
$ws.Range("E2902").Value = "Elektrik"
$ws.Range("F2902").Value = "Avtomat sığortalar və qutuları"
$ws.Range("G2902").Value = "Andeli Sığortalar"
$ws.Range("C2903").Value = 40.05
$ws.Range("D2903").Value = "Elektrik"
$ws.Range("E2903").Value = "Elektrik"
$ws.Range("F2903").Value = "Avtomat sığortalar və qutuları"
$ws.Range("G2903").Value = "Andeli Sığortalar"
$ws.Range("C2904").Value = 50.19
$ws.Range("D2904").Value = "Elektrik"
$ws.Range("E2904").Value = "Elektrik"
$ws.Range("F2904").Value = "Avtomat sığortalar və qutuları"
$ws.Range("G2904").Value = "Andeli Sığortalar"
$ws.Range("C2905").Value = 58.71
$ws.Range("D2905").Value = "Elektrik"
$ws.Range("E2905").Value = "Elektrik"
$ws.Range("F2905").Value = "Avtomat sığortalar və qutuları"
$ws.Range("G2905").Value = "Andeli Sığortalar"
$ws.Range("C2906").Value = 112.6
$ws.Range("D2906").Value = "Elektrik"
$ws.Range("E2906").Value = "Elektrik"
$ws.Range("F2906").Value = "Avtomat sığortalar və qutuları"
$ws.Range("G2906").Value = "Andeli Sığortalar"
$ws.Range("C2907").Value = 625
$ws.Range("D2907").Value = "Elektrik"
$ws.Range("E2907").Value = "Elektrikli Əl Alətləri"
$ws.Range("F2907").Value = "Qaynaq Aparatları və Aksesuarları"
$ws.Range("G2907").Value = "Metal Qaynaq Aparatları və Aksesuarları"
$ws.Range("C2908").Value = 643
$ws.Range("D2908").Value = "Elektrik"
$ws.Range("E2908").Value = "Elektrikli Əl Alətləri"
$ws.Range("F2908").Value = "Qaynaq Aparatları və Aksesuarları"
$ws.Range("G2908").Value = "Metal Qaynaq Aparatları və Aksesuarları"
$ws.Range("C2909").Value = 986
$ws.Range("D2909").Value = "Elektrik"
$ws.Range("E2909").Value = "Elektrikli Əl Alətləri"
$ws.Range("F2909").Value = "Qaynaq Aparatları və Aksesuarları"
$ws.Range("G2909").Value = "Metal Qaynaq Aparatları və Aksesuarları"
$ws.Range("C2910").Value = 240
$ws.Range("D2910").Value = "Elektrik"
$ws.Range("E2910").Value = "Elektrikli Əl Alətləri"
$ws.Range("F2910").Value = "Qaynaq Aparatları və Aksesuarları"
$ws.Range("G2910").Value = "Metal Qaynaq Aparatları və Aksesuarları"
$ws.Range("C2911").Value = 291
$ws.Range("D2911").Value = "Elektrik"
$ws.Range("E2911").Value = "Elektrikli Əl Alətləri"
$ws.Range("F2911").Value = "Qaynaq Aparatları və Aksesuarları"
$ws.Range("G2911").Value = "Metal Qaynaq Aparatları və Aksesuarları"
$ws.Range("C2912").Value = 320
$ws.Range("D2912").Value = "Elektrik"
$ws.Range("E2912").Value = "Elektrikli Əl Alətləri"
$ws.Range("F2912").Value = "Qaynaq Aparatları və Aksesuarları"
$ws.Range("G2912").Value = "Metal Qaynaq Aparatları və Aksesuarları"
$ws.Range("C2913").Value = 36
$ws.Range("D2913").Value = "Elektrik"
$ws.Range("E2913").Value = "Elektrik"
$ws.Range("F2913").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2913").Value = "Stablizatorlar"
$ws.Range("C2914").Value = 18.6
$ws.Range("D2914").Value = "Elektrik"
$ws.Range("E2914").Value = "Elektrik"
$ws.Range("F2914").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2914").Value = "Stablizatorlar"
$ws.Range("C2915").Value = 15.9
$ws.Range("D2915").Value = "Elektrik"
$ws.Range("E2915").Value = "Elektrik"
$ws.Range("F2915").Value = "Avtomat sığortalar və qutuları"
$ws.Range("G2915").Value = "Puskatel və Relelər"
$ws.Range("C2916").Value = 1789
$ws.Range("D2916").Value = "Elektrik"
$ws.Range("E2916").Value = "Elektrik"
$ws.Range("F2916").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2916").Value = "Stablizatorlar"
$ws.Range("C2917").Value = 1979
$ws.Range("D2917").Value = "Elektrik"
$ws.Range("E2917").Value = "Elektrik"
$ws.Range("F2917").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2917").Value = "Stablizatorlar"
$ws.Range("C2918").Value = 35
$ws.Range("D2918").Value = "Elektrik"
$ws.Range("E2918").Value = "Elektrik"
$ws.Range("F2918").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2918").Value = "Stablizatorlar"
$ws.Range("C2919").Value = 369
$ws.Range("D2919").Value = "Elektrik"
$ws.Range("E2919").Value = "Elektrikli Əl Alətləri"
$ws.Range("F2919").Value = "Qaynaq Aparatları və Aksesuarları"
$ws.Range("G2919").Value = "Metal Qaynaq Aparatları və Aksesuarları"
$ws.Range("C2920").Value = 329
$ws.Range("D2920").Value = "Elektrik"
$ws.Range("E2920").Value = "Elektrikli Əl Alətləri"
$ws.Range("F2920").Value = "Qaynaq Aparatları və Aksesuarları"
$ws.Range("G2920").Value = "Metal Qaynaq Aparatları və Aksesuarları"
$ws.Range("C2921").Value = 159
$ws.Range("D2921").Value = "Elektrik"
$ws.Range("E2921").Value = "Elektrik"
$ws.Range("F2921").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2921").Value = "Stablizatorlar"
$ws.Range("C2922").Value = 209
$ws.Range("D2922").Value = "Elektrik"
$ws.Range("E2922").Value = "Elektrik"
$ws.Range("F2922").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2922").Value = "Stablizatorlar"
$ws.Range("C2923").Value = 265
$ws.Range("D2923").Value = "Elektrik"
$ws.Range("E2923").Value = "Elektrik"
$ws.Range("F2923").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2923").Value = "Stablizatorlar"
$ws.Range("C2924").Value = 359
$ws.Range("D2924").Value = "Elektrik"
$ws.Range("E2924").Value = "Elektrik"
$ws.Range("F2924").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2924").Value = "Stablizatorlar"
$ws.Range("C2925").Value = 615
$ws.Range("D2925").Value = "Elektrik"
$ws.Range("E2925").Value = "Elektrik"
$ws.Range("F2925").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2925").Value = "Stablizatorlar"
$ws.Range("C2926").Value = 759
$ws.Range("D2926").Value = "Elektrik"
$ws.Range("E2926").Value = "Elektrik"
$ws.Range("F2926").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2926").Value = "Stablizatorlar"
$ws.Range("C2927").Value = 114
$ws.Range("D2927").Value = "Elektrik"
$ws.Range("E2927").Value = "Elektrik"
$ws.Range("F2927").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2927").Value = "Stablizatorlar"
$ws.Range("C2928").Value = 59.7
$ws.Range("D2928").Value = "Elektrik"
$ws.Range("E2928").Value = "Elektrik"
$ws.Range("F2928").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2928").Value = "Stablizatorlar"
$ws.Range("C2929").Value = 83.9
$ws.Range("D2929").Value = "Elektrik"
$ws.Range("E2929").Value = "Elektrik"
$ws.Range("F2929").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2929").Value = "Stablizatorlar"
$ws.Range("C2930").Value = 104
$ws.Range("D2930").Value = "Elektrik"
$ws.Range("E2930").Value = "Elektrik"
$ws.Range("F2930").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2930").Value = "Stablizatorlar"
$ws.Range("C2931").Value = 179
$ws.Range("D2931").Value = "Elektrik"
$ws.Range("E2931").Value = "Elektrik"
$ws.Range("F2931").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2931").Value = "Stablizatorlar"
$ws.Range("C2932").Value = 219
$ws.Range("D2932").Value = "Elektrik"
$ws.Range("E2932").Value = "Elektrik"
$ws.Range("F2932").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2932").Value = "Stablizatorlar"
$ws.Range("C2933").Value = 279
$ws.Range("D2933").Value = "Elektrik"
$ws.Range("E2933").Value = "Elektrik"
$ws.Range("F2933").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2933").Value = "Stablizatorlar"
$ws.Range("C2934").Value = 389
$ws.Range("D2934").Value = "Elektrik"
$ws.Range("E2934").Value = "Elektrik"
$ws.Range("F2934").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2934").Value = "Stablizatorlar"
$ws.Range("C2935").Value = 799
$ws.Range("D2935").Value = "Elektrik"
$ws.Range("E2935").Value = "Elektrik"
$ws.Range("F2935").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2935").Value = "Stablizatorlar"
$ws.Range("C2936").Value = 979
$ws.Range("D2936").Value = "Elektrik"
$ws.Range("E2936").Value = "Elektrik"
$ws.Range("F2936").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2936").Value = "Stablizatorlar"
$ws.Range("C2937").Value = 1149
$ws.Range("D2937").Value = "Elektrik"
$ws.Range("E2937").Value = "Elektrik"
$ws.Range("F2937").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2937").Value = "Stablizatorlar"
$ws.Range("C2938").Value = 2099
$ws.Range("D2938").Value = "Elektrik"
$ws.Range("E2938").Value = "Elektrik"
$ws.Range("F2938").Value = "Elektrik tənzimləyiciləri"
$ws.Range("G2938").Value = "Stablizatorlar"
$ws.Range("C2939").Value = 25.08
$ws.Range("D2939").Value = "Elektrik"
$ws.Range("E2939").Value = "Elektrik"
$ws.Range("F2939").Value = "Avtomat sığortalar və qutuları"
$ws.Range("G2939").Value = "Andeli Sığortalar"
$ws.Range("C2940").Value = 25.08
$ws.Range("D2940").Value = "Elektrik"
$ws.Range("E2940").Value = "Elektrik"
$ws.Range("F2940").Value = "Avtomat sığortalar və qutuları"
$ws.Range("G2940").Value = "Andeli Sığortalar"
$ws.Range("C2941").Value = 31.2
$ws.Range("D2941").Value = "Elektrik"
$ws.Range("E2941").Value = "Elektrik"
$ws.Range("F2941").Value = "Avtomat sığortalar və qutuları"
$ws.Range("G2941").Value = "Andeli Sığortalar"
$ws.Range("C2942").Value = 45.03
$ws.Range("D2942").Value = "Elektrik"
$ws.Range("E2942").Value = "Elektrik"
$ws.Range("F2942").Value = "Avtomat sığortalar və qutuları"
$ws.Range("G2942").Value = "Andeli Sığortalar"
$ws.Range("C2943").Value = 45.03
$ws.Range("D2943").Value = "Elektrik"
$ws.Range("E2943").Value = "Elektrik"
$ws.Range("F2943").Value = "Avtomat sığortalar və qutuları"
$ws.Range("G2943").Value = "Andeli Sığortalar"

# Phase 3: brand (H) for all new rows, in row order
$ws.Range("H2892").Value = "ANDELI"
$ws.Range("H2893").Value = "ANDELI"
$ws.Range("H2894").Value = "ANDELI"
$ws.Range("H2895").Value = "ANDELI"
$ws.Range("H2896").Value = "ANDELI"
$ws.Range("H2897").Value = "ANDELI"
$ws.Range("H2898").Value = "ANDELI"
$ws.Range("H2899").Value = "ANDELI"
$ws.Range("H2900").Value = "ANDELI"
$ws.Range("H2901").Value = "ANDELI"
$ws.Range("H2902").Value = "ANDELI"
$ws.Range("H2903").Value = "ANDELI"
$ws.Range("H2904").Value = "ANDELI"
$ws.Range("H2905").Value = "ANDELI"
$ws.Range("H2906").Value = "ANDELI"
$ws.Range("H2907").Value = "ANDELI"
$ws.Range("H2908").Value = "ANDELI"
$ws.Range("H2909").Value = "ANDELI"
$ws.Range("H2910").Value = "ANDELI"
$ws.Range("H2911").Value = "ANDELI"
$ws.Range("H2912").Value = "ANDELI"
$ws.Range("H2913").Value = "ANDELI"
$ws.Range("H2914").Value = "ANDELI"
$ws.Range("H2915").Value = "ANDELI"
$ws.Range("H2916").Value = "ANDELI"
$ws.Range("H2917").Value = "ANDELI"
$ws.Range("H2918").Value = "ANDELI"
$ws.Range("H2919").Value = "ANDELI"
$ws.Range("H2920").Value = "ANDELI"
$ws.Range("H2921").Value = "ANDELI"
$ws.Range("H2922").Value = "ANDELI"
$ws.Range("H2923").Value = "ANDELI"
$ws.Range("H2924").Value = "ANDELI"
$ws.Range("H2925").Value = "ANDELI"
$ws.Range("H2926").Value = "ANDELI"
$ws.Range("H2927").Value = "ANDELI"
$ws.Range("H2928").Value = "ANDELI"
$ws.Range("H2929").Value = "ANDELI"
$ws.Range("H2930").Value = "ANDELI"
$ws.Range("H2931").Value = "ANDELI"
$ws.Range("H2932").Value = "ANDELI"
$ws.Range("H2933").Value = "ANDELI"
$ws.Range("H2934").Value = "ANDELI"
$ws.Range("H2935").Value = "ANDELI"
$ws.Range("H2936").Value = "ANDELI"
$ws.Range("H2937").Value = "ANDELI"
$ws.Range("H2938").Value = "ANDELI"
$ws.Range("H2939").Value = "ANDELI"
$ws.Range("H2940").Value = "ANDELI"
$ws.Range("H2941").Value = "ANDELI"
$ws.Range("H2942").Value = "ANDELI"
$ws.Range("H2943").Value = "ANDELI"

# Update the active window view to match the post-edit state:
# top-left visible cell at A2891 and active selection at D2892
$win = $excel.ActiveWindow
$win.ScrollRow = 2891
$win.ScrollColumn = 1
$ws.Range("D2892").Select()
